$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A3"   = -22.13350000000001
    "D18"  = -8.475700000000002
    "A21"  = -20.16339999999998
    "A23"  = -20.26899999999998
    "C24"  = -13.33639999999999
    "A25"  = -21.67429999999998
    "C28"  = -13.29499999999999
    "C36"  = -11.78460000000001
    "C45"  = -14.13369999999999
    "C48"  = -11.5862
    "C49"  = -13.50039999999999
    "D51"  = -8.215500000000006
    "C52"  = -10.8827
    "A53"  = -22.1332
    "C53"  = -12.6853
    "C54"  = -13.7321
    "D55"  = -8.6097
    "A57"  = -22.58870000000002
    "A59"  = -22.0914
    "D64"  = -7.807099999999991
    "A69"  = -21.63100000000002
    "C70"  = -11.6473
    "A79"  = -20.46850000000001
    "D80"  = -7.757399999999997
    "A83"  = -21.87930000000001
    "C86"  = -14.26699999999999
    "C87"  = -13.4295
    "D92"  = -6.4623
    "A93"  = -21.43910000000001
    "D94"  = -6.557700000000002
    "D96"  = -8.575299999999997
    "C101" = -13.21469999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
